# TD_PPUIValidations.xlsx - "Add files via upload"
#
# The FunderPaid scenario column (column B) on the UIValidations sheet still
# priced everything in EUR while the Society/Multiple scenario columns
# (C and D) had already moved to USD. Align column B with the other two
# scenario columns by switching its EUR money cells to the equivalent USD
# figures that are already used elsewhere on the sheet (rows 28-33: the
# BaseAPCPrice/BaseArticleTypeDiscount/BaseAPCCharge/FinalNetPrice/Tax/
# TotalCharge block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UIValidations")

$ws.Range("B28").Value = "1,500.00 USD"   # BaseAPCPrice        (was 1,267.00 EUR)
$ws.Range("B29").Value = "0.00 USD"       # BaseArticleTypeDiscount (was 0.00 EUR)
$ws.Range("B30").Value = "1,500.00 USD"   # BaseAPCCharge       (was 1,267.00 EUR)
$ws.Range("B31").Value = "1,500.00 USD"   # FinalNetPrice       (was 1,267.00 EUR)
$ws.Range("B32").Value = "0.00 USD"       # Tax                 (was 0.00 EUR)
$ws.Range("B33").Value = "1,500.00 USD"   # TotalCharge         (was 1,267.00 EUR)

# Match the author's on-screen state: cell D31 selected, scrolled down so
# row 22 is visible at the top of the window.
$excel.Goto($ws.Range("A22"), $true) | Out-Null
$ws.Range("D31").Select() | Out-Null
